# ---------------------------------------------------------------------------
# Applies the 2025-05-09 trend-summary refresh:
#   * "Summary Table"  - drop the placeholder "------" divider row, trim the
#                          header captions, restyle the header row, and reset
#                          the page margins to the workbook defaults.
#   * "Cooccurrence"   - populate the (previously empty) keyword co-occurrence
#                          table with a header row + two data rows.
#   * "Associations"   - replace the placeholder divider row with the real
#                          first association and append the remaining rows.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Summary Table" sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary Table")

# Remove the "---------" placeholder divider that used to sit in row 2; this
# shifts every data row up by one (row 3 -> row 2 ... row 23 -> row 22) while
# preserving their contents untouched.
$summary.Rows.Item(2).Delete()

# Re-write the header captions without the leading/trailing padding spaces.
$summary.Range("A1").Value2 = "Keyword"
$summary.Range("B1").Value2 = "Keyword Count"
$summary.Range("C1").Value2 = "Short Summary"
$summary.Range("D1").Value2 = "Source URL"
$summary.Range("E1").Value2 = "Detailed Summary"

# Give the header the same bold / bordered / centered look used by the other
# generated sheets (Cooccurrence, Associations).
$summaryHeader = $summary.Range("A1:E1")
$summaryHeader.Font.Bold = $true
$summaryHeader.HorizontalAlignment = -4108   # xlCenter
$summaryHeader.VerticalAlignment = -4160     # xlTop
$summaryHeader.Borders.LineStyle = 1         # xlContinuous

# Reset page margins back to Excel's stock defaults (0.75/0.75/1/1/0.5/0.5 in).
$summary.PageSetup.LeftMargin = 54
$summary.PageSetup.RightMargin = 54
$summary.PageSetup.TopMargin = 72
$summary.PageSetup.BottomMargin = 72
$summary.PageSetup.HeaderMargin = 36
$summary.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# 2) "Cooccurrence" sheet - was completely empty, now gets a small table.
# ---------------------------------------------------------------------------
$cooc = $wb.Worksheets.Item("Cooccurrence")

$cooc.Cells.Item(1, 1).Value2 = "source"
$cooc.Cells.Item(1, 2).Value2 = "target"
$cooc.Cells.Item(1, 3).Value2 = "count"

$coocHeader = $cooc.Range("A1:C1")
$coocHeader.Font.Bold = $true
$coocHeader.HorizontalAlignment = -4108
$coocHeader.VerticalAlignment = -4160
$coocHeader.Borders.LineStyle = 1

$cooc.Cells.Item(2, 1).Value2 = "人工智能"
$cooc.Cells.Item(2, 2).Value2 = "科技成果转化"
$cooc.Cells.Item(2, 3).Value2 = 2

$cooc.Cells.Item(3, 1).Value2 = "人工智能"
$cooc.Cells.Item(3, 2).Value2 = "技术标准互认"
$cooc.Cells.Item(3, 3).Value2 = 1

# ---------------------------------------------------------------------------
# 3) "Associations" sheet - replace the placeholder row and append the rest.
# ---------------------------------------------------------------------------
$assoc = $wb.Worksheets.Item("Associations")

$assoc.Cells.Item(2, 1).Value2 = "新质生产力"
$assoc.Cells.Item(2, 2).Value2 = 2

$assoc.Cells.Item(3, 1).Value2 = "人工智能"
$assoc.Cells.Item(3, 2).Value2 = 4

$assoc.Cells.Item(4, 1).Value2 = "科技成果转化"
$assoc.Cells.Item(4, 2).Value2 = 3

$assoc.Cells.Item(5, 1).Value2 = "量子通信"
$assoc.Cells.Item(5, 2).Value2 = 1

$assoc.Cells.Item(6, 1).Value2 = "量子计算云平台"
$assoc.Cells.Item(6, 2).Value2 = 1

$assoc.Cells.Item(7, 1).Value2 = "生物技术"
$assoc.Cells.Item(7, 2).Value2 = 1

$assoc.Cells.Item(8, 1).Value2 = "碳纤维复合材料"
$assoc.Cells.Item(8, 2).Value2 = 1

$assoc.Cells.Item(9, 1).Value2 = "工业互联网安全"
$assoc.Cells.Item(9, 2).Value2 = 1

$assoc.Cells.Item(10, 1).Value2 = "技术标准互认"
$assoc.Cells.Item(10, 2).Value2 = 1

$assoc.Cells.Item(11, 1).Value2 = "智慧城市"
$assoc.Cells.Item(11, 2).Value2 = 1

# Leave the originally-active sheet selected.
$summary.Activate()
